# Daily scrape update - 2025-08-26 03:15:28 UTC
# Updates the opportunities table: new rows of scraped data replace the old
# ones, the table shrinks from 10 data rows to 5, a couple of column widths
# are narrowed, and the new "PREMIUM = Yes" cell on row 2 gets a yellow
# highlight fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Overwrite data rows 2-6 with the freshly scraped opportunities.
# ---------------------------------------------------------------------

$ws.Range("A2").Value = "1314884"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1314884"
$ws.Range("C2").Value = "ACE Program | Accounts Receivable Specialist"
$ws.Range("D2").Value = "Mumbai, Maharashtra, India"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "62 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "Tata Consultancy Services Ltd."

$ws.Range("A3").Value = "1327081"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327081"
$ws.Range("C3").Value = "Cyber Security Intern"
$ws.Range("D3").Value = "Chandigarh, India"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "2 applicants"
$ws.Range("G3").Value = "3 - 6 Months"
$ws.Range("H3").Value = "Solitaire Infosys Pvt. Ltd"

$ws.Range("A4").Value = "1327080"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1327080"
$ws.Range("C4").Value = "Machine Learning Intern"
$ws.Range("D4").Value = "Chandigarh, India"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "1 applicant"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "Solitaire Infosys Pvt. Ltd"

$ws.Range("A5").Value = "1327079"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1327079"
$ws.Range("C5").Value = "Software Developer Intern"
$ws.Range("D5").Value = "Chandigarh, India"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "3 applicants"
$ws.Range("G5").Value = "3 - 6 Months"
$ws.Range("H5").Value = "Solitaire Infosys Pvt. Ltd"

$ws.Range("A6").Value = "1325986"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1325986"
$ws.Range("C6").Value = "[Partly Remote] Internal Control Analyst"
$ws.Range("D6").Value = "Mexico City, CDMX, Mexico"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "33 applicants"
$ws.Range("G6").Value = "Partly Remote"
$ws.Range("H6").Value = "Sodexo Mexico"

# ---------------------------------------------------------------------
# 2. The scrape now only yields 5 listings instead of 10 - drop the
#    now-unused trailing rows so the sheet's used range becomes A1:H6.
# ---------------------------------------------------------------------

$ws.Rows("7:11").Delete()

# ---------------------------------------------------------------------
# 3. Highlight the new premium opportunity (row 2) with a yellow fill.
# ---------------------------------------------------------------------

$ws.Range("E2").Interior.Color = 65535

# ---------------------------------------------------------------------
# 4. Narrow a few columns to better fit the refreshed content.
# ---------------------------------------------------------------------

$ws.Columns(3).ColumnWidth = 46.1666667
$ws.Columns(4).ColumnWidth = 28.1666667
$ws.Columns(8).ColumnWidth = 32.1666667
